$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2033919902778603
$ws.Range("C2").Value = 0.4364632531330045
$ws.Range("D2").Value = 0.3742668131629021
$ws.Range("E2").Value = 0.6117734982515196
$ws.Range("F2").Value = 0.5987537801616326

$ws.Range("B3").Value = 0.2628942463139419
$ws.Range("C3").Value = 0.4333144524565598
$ws.Range("D3").Value = 0.2926323150708423
$ws.Range("E3").Value = 0.5409550028152456
$ws.Range("F3").Value = 0.4920828262122367

$ws.Range("B4").Value = 0.2853914118626144
$ws.Range("C4").Value = 0.4746291375408727
$ws.Range("D4").Value = 0.3229524026090269
$ws.Range("E4").Value = 0.5682890132749593
$ws.Range("F4").Value = 0.5132826384015559

$ws.Range("B5").Value = 0.2197630776494257
$ws.Range("C5").Value = 0.4033531803712247
$ws.Range("D5").Value = 0.234430775428361
$ws.Range("E5").Value = 0.4841805194639299
$ws.Range("F5").Value = 0.4524913939993277

$ws.Range("B6").Value = 0.1926414270218137
$ws.Range("C6").Value = 0.3605436531697729
$ws.Range("D6").Value = 0.1752150039837381
$ws.Range("E6").Value = 0.4185869132972722
$ws.Range("F6").Value = 0.3917259310889259

$ws.Range("B7").Value = 0.2387356972361521
$ws.Range("C7").Value = 0.2818843305516119
$ws.Range("D7").Value = 0.1725877670096401
$ws.Range("E7").Value = 0.4154368387729236
$ws.Range("F7").Value = 0.3606135925185842

$ws.Range("B8").Value = 0.3398310944915218
$ws.Range("C8").Value = 0.3398310944915218
$ws.Range("D8").Value = 0.2116720379002942
$ws.Range("E8").Value = 0.4600782954023089
$ws.Range("F8").Value = 0.3397414283545448

$ws.Range("B9").Value = 0.4537447083867415
$ws.Range("C9").Value = 0.4537447083867415
$ws.Range("D9").Value = 0.2538588652253824
$ws.Range("E9").Value = 0.5038440882112069
$ws.Range("F9").Value = 0.2682571662688996

$ws.Range("B10").Value = 0.08161590864515353
$ws.Range("C10").Value = 0.08161590864515353
$ws.Range("D10").Value = 0.006661156543974047
$ws.Range("E10").Value = 0.08161590864515353
